$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 (05-19-2015 / Tuesday) gets filled in with Official Business data,
#     matching the same visual style (fill/border/font) used by the
#     surrounding populated rows (e.g. row 16). Copy formats only, so the
#     existing values in row 15 are preserved.
$srcRow = $ws.Range("A16:P16")
$dstRow = $ws.Range("A15:P15")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)

# Fill in the Official Business columns for row 15.
$ws.Range("K15").Value = "06:30:00"
$ws.Range("L15").Value = "08:30:00"
$ws.Range("M15").Value = "16:30:00"
$ws.Range("N15").Value = "18:30:00"
$ws.Range("P15").Value = "~OB Others|Integration test| R"

# --- Append " R" to the remarks that reference the shared "Integration
#     test(ing)" strings (this also updates row 16's identical remark,
#     since it shares the same underlying text).
$ws.Range("P8").Value = "~OB Others|Integration Testing| R"
$ws.Range("P9").Value = "~OB Others|Integration Test| R"
$ws.Range("P10").Value = "~OB Others|Integration test| R"
$ws.Range("P16").Value = "~OB Others|Integration test| R"
$ws.Range("P14").Value = "~OT ~ = Adjusted processBonusBuy of barterSapService to handle a scenario of Buy 2 of any subItems for a fixed price ~OB Others|integration test| R"

# --- Total overtime hours adjustment.
$ws.Range("I23").Value = 2
